$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Generators")

# Replace the hard-coded values in column H (rows 2-54) with a formula that
# mirrors column G with the opposite sign: H = -G
$ws.Range("H2").Formula = "=-G2"
$ws.Range("H3:H54").Formula = "=-G3"

# Reflect the selection left behind on the Generators sheet after the edit
$ws.Range("H2:H54").Select()
